$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-19 changes from serial date 45179 (2023-09-10)
# to serial date 45180 (2023-09-11).
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
